$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.162.87"
$ws.Range("E2").Value = "  -8.43%  "
$ws.Range("D3").Value = "3.177.87"
$ws.Range("E3").Value = "  -10.02%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "506.11"
$ws.Range("E5").Value = "  -9.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.18"
$ws.Range("E6").Value = "  -15.73%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.583"
$ws.Range("E7").Value = "  -8.16%  "
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D9").Value = "3.175.44"
$ws.Range("E9").Value = "  -9.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.584"
$ws.Range("E10").Value = "  -11.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.66"
$ws.Range("E11").Value = "  -13.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.128"
$ws.Range("E12").Value = "  -11.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000247"
$ws.Range("E13").Value = "  -9.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.70"
$ws.Range("E14").Value = "  -12.83%  "
$ws.Range("D15").Value = "3.701.07"
$ws.Range("E15").Value = "  -9.40%  "
$ws.Range("D16").Value = "3.194.51"
$ws.Range("E16").Value = "  -9.29%  "
$ws.Range("D17").Value = "62.167.83"
$ws.Range("E17").Value = "  -8.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.111"
$ws.Range("E18").Value = "  -10.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.74"
$ws.Range("E19").Value = "  -9.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.59"
$ws.Range("E20").Value = "  -11.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.923"
$ws.Range("E21").Value = "  -10.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "355.79"
$ws.Range("E22").Value = "  -10.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.61"
$ws.Range("E23").Value = "  -9.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.83"
$ws.Range("E24").Value = "  -8.75%  "
$ws.Range("B25").Value = "LEO"
$ws.Range("C25").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.06"
$ws.Range("E25").Value = "  -1.84%  "
$ws.Range("B26").Value = "RenderToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.51"
$ws.Range("E26").Value = "  -11.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.78"
$ws.Range("E27").Value = "  -2.67%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.56"
$ws.Range("E28").Value = "  -9.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.86"
$ws.Range("E29").Value = "  -12.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.96"
$ws.Range("E30").Value = "  -10.78%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "27.53"
$ws.Range("E31").Value = "  -12.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "609.81"
$ws.Range("E32").Value = "  -15.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.31"
$ws.Range("E33").Value = "  -11.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.84"
$ws.Range("E34").Value = "  -7.90%  "
$ws.Range("E35").Value = "  -9.99%  "
$ws.Range("B36").Value = "Dai"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "56.09"
$ws.Range("E37").Value = "  -13.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "35.41"
$ws.Range("E38").Value = "  -8.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.368"
$ws.Range("E39").Value = "  -7.16%  "
$ws.Range("E40").Value = "  +0.17%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "2.813.35"
$ws.Range("E41").Value = "  -8.78%  "
$ws.Range("B42").Value = "PEPE"
$ws.Range("C42").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D42").Value = "0.0₃0654"
$ws.Range("E42").Value = "  -4.76%  "
$ws.Range("E43").Value = "  -10.90%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.34"
$ws.Range("E44").Value = "  -6.78%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.56"
$ws.Range("E45").Value = "  -7.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.55"
$ws.Range("E46").Value = "  -15.94%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0377"
$ws.Range("E47").Value = "  -7.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.90"
$ws.Range("E48").Value = "  -1.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.67"
$ws.Range("E49").Value = "  +0.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.30"
$ws.Range("E50").Value = "  -5.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.120"
$ws.Range("E51").Value = "  -8.81%  "
